$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from the existing "sum" header cell (G1) onto the
# new "Save" header cell (H1) so it picks up the same bold/border/centered
# formatting as the rest of row 1, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New data value for the "Save" column.
$ws.Range("H2").Value = 0
